$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.14588451084694
$ws.Range("C2").Value = 0.2026105229323036
$ws.Range("E2").Value = 0.09905265697044818
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.002427098136627533
$ws.Range("L2").Value = 0.2110524580248097
$ws.Range("N2").Value = 1.16840010077054
$ws.Range("O2").Value = 2.566964032274484

$ws.Range("B3").Value = 1.048443709493881
$ws.Range("C3").Value = 0.1940283572276797
$ws.Range("E3").Value = 0.09955546178018793
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.002429989401759624
$ws.Range("L3").Value = 0.201969337403483
$ws.Range("N3").Value = 1.177748425787925
$ws.Range("O3").Value = 2.572591318806104

$ws.Range("B4").Value = 0.9888365736096034
$ws.Range("C4").Value = 0.1887215868100611
$ws.Range("E4").Value = 0.09991664485843899
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.002431859756739937
$ws.Range("L4").Value = 0.1964989048526462
$ws.Range("N4").Value = 1.183992298864418
$ws.Range("O4").Value = 2.577941223246484

$ws.Range("B5").Value = 0.9646032566734277
$ws.Range("C5").Value = 0.1865497991471017
$ws.Range("E5").Value = 0.1000770239804947
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002432645928773595
$ws.Range("L5").Value = 0.1942965388698497
$ws.Range("N5").Value = 1.186663572456069
$ws.Range("O5").Value = 2.5805969591203

$ws.Range("B6").Value = 0.9605828176650562
$ws.Range("C6").Value = 0.1861886217879061
$ws.Range("E6").Value = 0.1001044517878462
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.002432777923163419
$ws.Range("L6").Value = 0.1939324630490802
$ws.Range("N6").Value = 1.187114800092019
$ws.Range("O6").Value = 2.581066642991601

$ws.Range("B7").Value = 0.9885095218378979
$ws.Range("C7").Value = 0.1886923345295912
$ws.Range("E7").Value = 0.09991875436096187
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.00243187026190649
$ws.Range("L7").Value = 0.1964690940534837
$ws.Range("N7").Value = 1.184027810846217
$ws.Range("O7").Value = 2.577975114833777

$ws.Range("B8").Value = 1.112241596862361
$ws.Range("C8").Value = 0.1996592124045407
$ws.Range("E8").Value = 0.09921513781474545
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.002428075348938673
$ws.Range("L8").Value = 0.2078984999562152
$ws.Range("N8").Value = 1.171518875393105
$ws.Range("O8").Value = 2.568510558605482

$ws.Range("B9").Value = 1.356597209175732
$ws.Range("C9").Value = 0.2208641878712569
$ws.Range("E9").Value = 0.09825152663344028
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.002421384862525773
$ws.Range("L9").Value = 0.2311562442706361
$ws.Range("N9").Value = 1.150982716287736
$ws.Range("O9").Value = 2.565025558152456

$ws.Range("B10").Value = 1.537132542428026
$ws.Range("C10").Value = 0.2362544131710536
$ws.Range("E10").Value = 0.09779729583559948
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002416922766212748
$ws.Range("L10").Value = 0.2487588058665295
$ws.Range("N10").Value = 1.138323425326973
$ws.Range("O10").Value = 2.57171700245641

$ws.Range("B11").Value = 1.619474234484073
$ws.Range("C11").Value = 0.2432136863365031
$ws.Range("E11").Value = 0.09764576324360874
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.002414990311196043
$ws.Range("L11").Value = 0.2568786761376884
$ws.Range("N11").Value = 1.133090574472334
$ws.Range("O11").Value = 2.576783689951299

$ws.Range("B12").Value = 1.650684823891481
$ws.Range("C12").Value = 0.2458428459636934
$ws.Range("E12").Value = 0.09759630471075198
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.002414272468622592
$ws.Range("L12").Value = 0.2599695866691434
$ws.Range("N12").Value = 1.131184574888074
$ws.Range("O12").Value = 2.578994196649887

$ws.Range("B13").Value = 1.643961771432544
$ws.Range("C13").Value = 0.2452768860604806
$ws.Range("E13").Value = 0.09760660408196387
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.002414426450172535
$ws.Range("L13").Value = 0.2593031891701685
$ws.Range("N13").Value = 1.131591706013481
$ws.Range("O13").Value = 2.578505127944567

$ws.Range("B14").Value = 1.62204136354859
$ws.Range("C14").Value = 0.2434301133316694
$ws.Range("E14").Value = 0.09764153545803111
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.002414930974906422
$ws.Range("L14").Value = 0.2571326450899392
$ws.Range("N14").Value = 1.132932252578591
$ws.Range("O14").Value = 2.576959694601214

$ws.Range("B15").Value = 1.608618291879679
$ws.Range("C15").Value = 0.2422981041821117
$ws.Range("E15").Value = 0.09766396382831033
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002415241823883205
$ws.Range("L15").Value = 0.255805217796464
$ws.Range("N15").Value = 1.133763215487043
$ws.Range("O15").Value = 2.576051110731328

$ws.Range("B16").Value = 1.531755533040894
$ws.Range("C16").Value = 0.2357987537523059
$ws.Range("E16").Value = 0.0978083076160523
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002417051010459593
$ws.Range("L16").Value = 0.2482304087616996
$ws.Range("N16").Value = 1.138675982144129
$ws.Range("O16").Value = 2.571426670047231

$ws.Range("B17").Value = 1.484656873015354
$ws.Range("C17").Value = 0.231800793391983
$ws.Range("E17").Value = 0.09791097035594554
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.002418185780582134
$ws.Range("L17").Value = 0.2436122429091512
$ws.Range("N17").Value = 1.141824462942559
$ws.Range("O17").Value = 2.569108501744893

$ws.Range("B18").Value = 1.457587313104398
$ws.Range("C18").Value = 0.2294973453644502
$ws.Range("E18").Value = 0.09797520543174087
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002418847639953535
$ws.Range("L18").Value = 0.2409665775339818
$ws.Range("N18").Value = 1.143684891517893
$ws.Range("O18").Value = 2.567965480577016

$ws.Range("B19").Value = 1.448425561310273
$ws.Range("C19").Value = 0.2287167676447268
$ws.Range("E19").Value = 0.09799784504564713
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.002419073310833896
$ws.Range("L19").Value = 0.2400726201401682
$ws.Range("N19").Value = 1.144323305228902
$ws.Range("O19").Value = 2.567611131850697

$ws.Range("B20").Value = 1.489668510277909
$ws.Range("C20").Value = 0.2322267906324669
$ws.Range("E20").Value = 0.0978995049790754
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.002418064033903877
$ws.Range("L20").Value = 0.2441027598969754
$ws.Range("N20").Value = 1.141484178580264
$ws.Range("O20").Value = 2.569335569255202

$ws.Range("B21").Value = 1.628479126923708
$ws.Range("C21").Value = 0.2439727238088096
$ws.Range("E21").Value = 0.09763106021971879
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.00241478240608912
$ws.Range("L21").Value = 0.2577697503115957
$ws.Range("N21").Value = 1.132536451128153
$ws.Range("O21").Value = 2.577405696829459

$ws.Range("B22").Value = 1.719371654871225
$ws.Range("C22").Value = 0.2516133590559662
$ws.Range("E22").Value = 0.09750180086199123
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.002412718873960884
$ws.Range("L22").Value = 0.2667956830351272
$ws.Range("N22").Value = 1.127129034973343
$ws.Range("O22").Value = 2.584381650431993

$ws.Range("B23").Value = 1.670845339304492
$ws.Range("C23").Value = 0.2475387535506854
$ws.Range("E23").Value = 0.09756656297840394
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.002413812811484805
$ws.Range("L23").Value = 0.2619698157860739
$ws.Range("N23").Value = 1.12997479276352
$ws.Range("O23").Value = 2.580502421130518

$ws.Range("B24").Value = 1.487402723775176
$ws.Range("C24").Value = 0.2320342127437129
$ws.Range("E24").Value = 0.09790467223860055
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.002418119046151787
$ws.Range("L24").Value = 0.2438809679469927
$ws.Range("N24").Value = 1.141637864372939
$ws.Range("O24").Value = 2.569232321120609

$ws.Range("B25").Value = 1.29031280641783
$ws.Range("C25").Value = 0.2151604845673347
$ws.Range("E25").Value = 0.09846765214371445
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.002423114861765207
$ws.Range("L25").Value = 0.2247739621767835
$ws.Range("N25").Value = 1.156111378306115
$ws.Range("O25").Value = 2.564348081102423
